$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new row (13) to the time-log sheet documenting the API pagination
# work done on 08.07.2019, mirroring the look & feel of the existing rows.
# ---------------------------------------------------------------------------

# 1) Copy the formatting of an existing "template" row (row 5) which already
#    has the exact style pattern we need for row 13:
#      - column A: default style (plain text date)
#      - column B: wrap-text style (style index 6)
#      - column C: default style (plain numeric hours)
$ws.Range("A5:C5").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# 2) Write the date into A13 as literal text "08.07.2019".
#    A plain .Value assignment of this particular string would be
#    auto-recognised as a real date (08 is a valid month number), so we
#    stage the text in a helper cell that has been explicitly marked as
#    Text, then copy only the *value* over (xlPasteValues) which preserves
#    the text typing without touching A13's formatting/style.
$helper = $ws.Cells.Item(1000, 1000)
$helper.NumberFormat = "@"
$helper.Value = "08.07.2019"
$helper.Copy()
$ws.Cells.Item(13, 1).PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
$helper.Clear()

# 3) Write the description text (column B) and the hours spent (column C).
$ws.Cells.Item(13, 2).Value = "The Card objects from the api are now returned `npaginated with additional info such as the link for the next page and the total number of cards and pages. This is because the json for all the cards is too big and slow."
$ws.Cells.Item(13, 3).Value = 1.5

# 4) Match the (word-wrapped, 4 visual lines) row height used by real Excel.
$ws.Rows.Item(13).RowHeight = 60

# 5) Update the active selection to reflect where the author finished editing.
$ws.Range("E13").Select()
